# Refresh leve-profit market figures (currentAveragePrice / LevePrice / LeveProfit
# columns) across several crafting-job sheets per the scheduled market-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3: One for the Books
$ws.Range("H3").Value = 67500
$ws.Range("J3").Value = 67500
$ws.Range("L3").Value = 67500
$ws.Range("N3").Value = -67728

# Row 6: Days of Chunder
$ws.Range("H6").Value = 3477.1667
$ws.Range("I6").Value = 4152.6
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 12457.8
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -12345.8
$ws.Range("N6").Value = -524

# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 298.43478
$ws.Range("I15").Value = 298.43478
$ws.Range("K15").Value = 895.3043399999999
$ws.Range("M15").Value = -726.3043399999999

# Row 87: There Was a Late Fee
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

# Row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

# Row 102: Spell-rebound
$ws.Range("H102").Value = 67500
$ws.Range("J102").Value = 67500
$ws.Range("L102").Value = 67500
$ws.Range("N102").Value = -73990

# Row 106: Making Your Mark
$ws.Range("H106").Value = 4032.6667
$ws.Range("I106").Value = 4032.6667
$ws.Range("K106").Value = 4032.6667
$ws.Range("M106").Value = -3401.6667

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 683.9
$ws.Range("I111").Value = 680
$ws.Range("K111").Value = 2040
$ws.Range("M111").Value = 1027

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3021.1177
$ws.Range("I132").Value = 3212.1428
$ws.Range("K132").Value = 9636.428400000001
$ws.Range("M132").Value = -7106.428400000001

$ws = $wb.Worksheets.Item("ARM")
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 2031
$ws.Range("I110").Value = 553.4
$ws.Range("K110").Value = 553.4
$ws.Range("M110").Value = 1491.6

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2084.4
$ws.Range("I132").Value = 2084.4
$ws.Range("K132").Value = 6253.200000000001
$ws.Range("M132").Value = -3723.200000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 3196.923
$ws.Range("I86").Value = 2313.6667
$ws.Range("J86").Value = 3954
$ws.Range("K86").Value = 2313.6667
$ws.Range("L86").Value = 3954
$ws.Range("M86").Value = -1190.6667
$ws.Range("N86").Value = -6200

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 3196.923
$ws.Range("I89").Value = 2313.6667
$ws.Range("J89").Value = 3954
$ws.Range("K89").Value = 11568.3335
$ws.Range("L89").Value = 19770
$ws.Range("M89").Value = -5952.333500000001
$ws.Range("N89").Value = -31002

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 1296.6666
$ws.Range("I105").Value = 1345
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 1345
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = 402
$ws.Range("N105").Value = -4694

$ws = $wb.Worksheets.Item("CRP")
# Row 26: As the Worm Turns
$ws.Range("H26").Value = 52000
$ws.Range("J26").Value = 52000
$ws.Range("L26").Value = 52000
$ws.Range("N26").Value = -52574

# Row 31: Wall Not Found
$ws.Range("H31").Value = 3231.3333
$ws.Range("I31").Value = 3231.3333
$ws.Range("K31").Value = 3231.3333
$ws.Range("M31").Value = -2936.3333

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3231.3333
$ws.Range("I34").Value = 3231.3333
$ws.Range("K34").Value = 3231.3333
$ws.Range("M34").Value = -3029.3333

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 3500
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 3500
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 99: O Pine
$ws.Range("H99").Value = 3004.625
$ws.Range("I99").Value = 3248.1428
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 3248.1428
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = -1750.1428
$ws.Range("N99").Value = -4296

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 2699.5
$ws.Range("I122").Value = 2600
$ws.Range("K122").Value = 7800
$ws.Range("M122").Value = -5350

# Row 126: A Better Conductor
$ws.Range("H126").Value = 3004.625
$ws.Range("I126").Value = 3248.1428
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 9744.428400000001
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -7274.428400000001
$ws.Range("N126").Value = -8840

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 10836.667
$ws.Range("I134").Value = 10004
$ws.Range("K134").Value = 30012
$ws.Range("M134").Value = -27477

$ws = $wb.Worksheets.Item("CUL")
# Row 37: I Love Lamprey
$ws.Range("H37").Value = 98420
$ws.Range("J37").Value = 98420
$ws.Range("L37").Value = 295260
$ws.Range("N37").Value = -295484

# Row 62: Little Orphan Candy
$ws.Range("H62").Value = 10999.8
$ws.Range("J62").Value = 11666.333
$ws.Range("L62").Value = 34998.999
$ws.Range("N62").Value = -36370.999

# Row 65: Confections of Confession (L)
$ws.Range("H65").Value = 10999.8
$ws.Range("J65").Value = 11666.333
$ws.Range("L65").Value = 104996.997
$ws.Range("N65").Value = -111860.997

# Row 122: Salt of the North
$ws.Range("H122").Value = 741.5714
$ws.Range("I122").Value = 781.8333
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 7036.4997
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -4586.4997
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 2670.3
$ws.Range("I7").Value = 2300.5
$ws.Range("J7").Value = 3225
$ws.Range("K7").Value = 2300.5
$ws.Range("L7").Value = 3225
$ws.Range("M7").Value = -2188.5
$ws.Range("N7").Value = -3449

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 2766.4443
$ws.Range("I46").Value = 2260.7273
$ws.Range("J46").Value = 3561.1428
$ws.Range("K46").Value = 2260.7273
$ws.Range("L46").Value = 3561.1428
$ws.Range("M46").Value = -2072.7273
$ws.Range("N46").Value = -3937.1428

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 4133.2104
$ws.Range("I61").Value = 3235.4666
$ws.Range("J61").Value = 7499.75
$ws.Range("K61").Value = 3235.4666
$ws.Range("L61").Value = 7499.75
$ws.Range("M61").Value = -3033.4666
$ws.Range("N61").Value = -7903.75

# Row 113: Peace in Rest
$ws.Range("H113").Value = 4133.2104
$ws.Range("I113").Value = 3235.4666
$ws.Range("J113").Value = 7499.75
$ws.Range("K113").Value = 3235.4666
$ws.Range("L113").Value = 7499.75
$ws.Range("M113").Value = -1065.4666
$ws.Range("N113").Value = -11839.75

# Row 126: Battered Books
$ws.Range("H126").Value = 2670.3
$ws.Range("I126").Value = 2300.5
$ws.Range("J126").Value = 3225
$ws.Range("K126").Value = 6901.5
$ws.Range("L126").Value = 9675
$ws.Range("M126").Value = -4431.5
$ws.Range("N126").Value = -14615

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax
$ws.Range("H107").Value = 478
$ws.Range("I107").Value = 516.3333
$ws.Range("K107").Value = 1548.9999
$ws.Range("M107").Value = 371.0001

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2825.125
$ws.Range("I126").Value = 2708
$ws.Range("J126").Value = 3332.6667
$ws.Range("K126").Value = 8124
$ws.Range("L126").Value = 9998.000100000001
$ws.Range("M126").Value = -5654
$ws.Range("N126").Value = -14938.0001
